$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 13.40353125
$ws.Range("Z2").Value = 11.6585
$ws.Range("AA2").Value = 20.97370052083333
$ws.Range("AB2").Value = 19.79076822916667

# Row 3
$ws.Range("D3").Value = 30.32099739583333
$ws.Range("E3").Value = 36.28490625
$ws.Range("F3").Value = 11.48641666666667
$ws.Range("N3").Value = 43.54557291666666
$ws.Range("O3").Value = 42.72395833333334
$ws.Range("P3").Value = 41.18515625
$ws.Range("Q3").Value = 55.234375
$ws.Range("R3").Value = 59.54557291666666
$ws.Range("S3").Value = 28.1703125
$ws.Range("Y3").Value = 6.787265625
$ws.Range("Z3").Value = 2.547158854166667
$ws.Range("AA3").Value = 2.1611484375
$ws.Range("AB3").Value = 4.280145833333332
$ws.Range("AD3").Value = 2.594625
$ws.Range("AE3").Value = 4.07925
$ws.Range("AF3").Value = 7.632166666666667
$ws.Range("AG3").Value = 10.18028385416667
$ws.Range("AH3").Value = 2.198145833333333

# Row 4
$ws.Range("D4").Value = 8.292645833333333
$ws.Range("F4").Value = 21.4059375
$ws.Range("H4").Value = 41.96613281250001
$ws.Range("X4").Value = 0.6979791666666666
$ws.Range("Y4").Value = 0.5233046875

# Row 5
$ws.Range("T5").Value = 0.5001302083333332
$ws.Range("U5").Value = 0.6860937499999999
$ws.Range("V5").Value = 0.09746093750000001
